$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 116; existing rows 116-149 shift down to 120-153
$ws.Rows.Item(116).Resize(4).Insert()

# Columns that stay constant across this whole price-list block
$constMercadoId = 11
$constMercado   = "Vega Monumental Concepción"
$constRegion    = "Bíobío"
$constCodreg    = 8
$constCatId     = 100112045
$constCategoria = "Zapallo"
$constUnidad    = "`$/kilo (volumen en unidades)"
$constKgUnid    = 1
$constClasif    = "Hortaliza"

for ($r = 116; $r -le 119; $r++) {
    $ws.Cells.Item($r, 1).Value  = $constMercadoId
    $ws.Cells.Item($r, 2).Value  = $constMercado
    $ws.Cells.Item($r, 3).Value  = $constRegion
    $ws.Cells.Item($r, 5).Value  = $constCodreg
    $ws.Cells.Item($r, 6).Value  = $constCatId
    $ws.Cells.Item($r, 7).Value  = $constCategoria
    $ws.Cells.Item($r, 14).Value = $constUnidad
    $ws.Cells.Item($r, 17).Value = $constKgUnid
    $ws.Cells.Item($r, 18).Value = $constClasif
}

# Row 116: Camote, 1a nueva(o)
$ws.Cells.Item(116, 4).Value  = 44559
$ws.Cells.Item(116, 8).Value  = "Camote"
$ws.Cells.Item(116, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(116, 10).Value = 300
$ws.Cells.Item(116, 11).Value = 550
$ws.Cells.Item(116, 12).Value = 550
$ws.Cells.Item(116, 13).Value = 550
$ws.Cells.Item(116, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(116, 16).Value = 550

# Row 117: Camote, 2a nueva(o)
$ws.Cells.Item(117, 4).Value  = 44559
$ws.Cells.Item(117, 8).Value  = "Camote"
$ws.Cells.Item(117, 9).Value  = "2a nueva(o)"
$ws.Cells.Item(117, 10).Value = 300
$ws.Cells.Item(117, 11).Value = 450
$ws.Cells.Item(117, 12).Value = 450
$ws.Cells.Item(117, 13).Value = 450
$ws.Cells.Item(117, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(117, 16).Value = 450

# Row 118: Paine, 1a nueva(o)
$ws.Cells.Item(118, 4).Value  = 44559
$ws.Cells.Item(118, 8).Value  = "Paine"
$ws.Cells.Item(118, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(118, 10).Value = 300
$ws.Cells.Item(118, 11).Value = 250
$ws.Cells.Item(118, 12).Value = 250
$ws.Cells.Item(118, 13).Value = 250
$ws.Cells.Item(118, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(118, 16).Value = 250

# Row 119: Paine, 2a nueva(o)
$ws.Cells.Item(119, 4).Value  = 44559
$ws.Cells.Item(119, 8).Value  = "Paine"
$ws.Cells.Item(119, 9).Value  = "2a nueva(o)"
$ws.Cells.Item(119, 10).Value = 300
$ws.Cells.Item(119, 11).Value = 200
$ws.Cells.Item(119, 12).Value = 200
$ws.Cells.Item(119, 13).Value = 200
$ws.Cells.Item(119, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(119, 16).Value = 200
